$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E14").Value = 101
$ws.Range("E15").Value = 100.49
$ws.Range("C19").Value = 98.56
$ws.Range("E19").Value = 96.69
$ws.Range("E20").Value = 99.33
$ws.Range("E22").Value = 98.25
